$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 493, shifting existing rows 493:548 down to 494:549
$ws.Rows("493:493").Insert()

# Populate the newly inserted row 493 with the new record
$ws.Range("A493").Value = 5
$ws.Range("B493").Value = "Macroferia Regional de Talca"
$ws.Range("C493").Value = "Maule"
$ws.Range("D493").Value = 44918
$ws.Range("E493").Value = 7
$ws.Range("F493").Value = 100112043
$ws.Range("G493").Value = "Pepino ensalada"
$ws.Range("H493").Value = "Sin especificar"
$ws.Range("I493").Value = "Primera"
$ws.Range("J493").Value = 500
$ws.Range("K493").Value = 12000
$ws.Range("L493").Value = 12000
$ws.Range("M493").Value = 12000
$ws.Range("N493").Value = "`$/caja 80 unidades"
$ws.Range("O493").Value = "Región del Maule"
$ws.Range("P493").Value = 150
$ws.Range("Q493").Value = 80
$ws.Range("R493").Value = "Hortaliza"
